$d = $word.ActiveDocument

$lines = @(
    '********************************8月28日读经章节***************************',
    'Chapter 1 of 1_Thessalonians',
    '1.保罗，西拉，提摩太，写信给帖撒罗尼迦在父神和主耶稣基督里的教会。愿恩惠平安归与你们。',
    '2.我们为你们众人常常感谢神，祷告的时候提到你们。',
    '3.在神我们的父面前，不住地记念你们因信心所作的工夫，因爱心所受的劳苦，因盼望我们主耶稣基督所存的忍耐。',
    '4.被神所爱的弟兄阿，我知道你们是蒙拣选的。',
    '5.因为我们的福音传到你们那里，不独在乎言语，也在乎权能和圣灵，并充足的信心，正如你们知道我们在你们那里，为你们的缘故是怎样为人。',
    '6.并且你们在大难之中，蒙了圣灵所赐的喜乐，领受真道，就效法我们，也效法了主。',
    '7.甚至你们作了马其顿和亚该亚，所有信主之人的榜样。',
    '8.因为主的道从你们那里已经传扬出来，你们向神的信心不但在马其顿和亚该亚，就是在各处，也都传开了。所以不用我们说什么话。',
    '9.因为他们自己已经报明我们是怎样进到你们那里，你们是怎样离弃偶像归向神，要服事那又真又活的神，',
    '10.等候他儿子从天降临，就是他从死里复活的，那位救我们脱离将来忿怒的耶稣。',
    'Chapter 16 of Proverbs',
    '1.心中的谋算在乎人。舌头的应对，由于耶和华。',
    '2.人一切所行的，在自己眼中看为清洁。惟有耶和华衡量人心。',
    '3.你所作的，要交托耶和华，你所谋的，就必成立。',
    '4.耶和华所造的，各适其用。就是恶人，也为祸患的日子所造。',
    '5.凡心里骄傲的，为耶和华所憎恶。虽然连手，他必不免受罚。',
    '6.因怜悯诚实，罪孽得赎。敬畏耶和华的，远离恶事。',
    '7.人所行的若蒙耶和华喜悦，耶和华也使他的仇敌与他和好。',
    '8.多有财利，行事不义，不如少有财利，行事公义。',
    '9.人心筹算自己的道路。惟耶和华指引他的脚步。',
    '10.王的嘴中有神语。审判之时，他的口，必不差错。',
    '11.公道的天平和秤，都属耶和华。囊中一切法码，都为他所定。',
    '12.作恶为王所憎恶。因国位是靠公义坚立。',
    '13.公义的嘴，为王所喜悦。说正直话的，为王所喜爱。',
    '14.王的震怒，如杀人的使者。但智慧人能止息王怒。',
    '15.王的脸光，使人有生命。王的恩典，好像春云时雨。',
    '16.得智慧胜似得金子。选聪明强如选银子。',
    '17.正直人的道，是远离恶事。谨守己路的，是保全性命。',
    '18.骄傲在败坏以先，狂心在跌倒之前。',
    '19.心里谦卑与穷乏人来往，强如将掳物与骄傲人同分。',
    '20.谨守训言的，必得好处。倚靠耶和华的，便为有福。',
    '21.心中有智慧，必称为通达人。嘴中的甜言，加增人的学问。',
    '22.人有智慧就有生命的泉源。愚昧人必被愚昧惩治。',
    '23.智慧人的心，教训他的口，又使他的嘴，增长学问。',
    '24.良言如同蜂房，使心觉甘甜，使骨得医治。',
    '25.有一条路，人以为正，至终成为死亡之路。',
    '26.劳力人的胃口，使他劳力，因为他的口腹催逼他。',
    '27.匪徒图谋奸恶，嘴上仿佛有烧焦的火。',
    '28.乖僻人播散分争。传舌的离间密友。',
    '29.强暴人诱惑邻舍，领他走不善之道。',
    '30.眼目紧合的，图谋乖僻，嘴唇紧闭的，成就邪恶。',
    '31.白发是荣耀的冠冕。在公义的道上，必能得着。',
    '32.不轻易发怒的，胜过勇士。治服己心的，强如取城。',
    '33.签放在怀里。定事由耶和华。',
    'Chapter 17 of Proverbs',
    '1.设筵满屋，大家相争，不如有块干饼，大家相安。',
    '2.仆人办事聪明，必管辖贻羞之子，又在众子中，同分产业。',
    '3.鼎为炼银，炉为炼金。惟有耶和华熬炼人心。',
    '4.行恶的留心听奸诈之言。说谎的侧耳听邪恶之语。',
    '5.戏笑穷人的，是辱没造他的主。幸灾乐祸的，必不免受罚。',
    '6.子孙为老人的冠冕。父亲是儿女的荣耀。',
    '7.愚顽人说美言本不相宜，何况君王说谎话呢？',
    '8.贿赂在馈送的人眼中，看为宝玉。随处运动，都得顺利。',
    '9.遮掩人过的，寻求人爱。屡次挑错的，离间密友。',
    '10.一句责备话，深入聪明人的心，强如责打愚昧人一百下。',
    '11.恶人只寻背叛，所以必有严厉的使者，奉差攻击他。',
    '12.宁可遇见丢崽子的母熊，不可遇见正行愚妄的愚昧人。',
    '13.以恶报善的，祸患必不离他的家。',
    '14.分争的起头，如水放开。所以在争闹之先，必当止息争竞。',
    '15.定恶人为义的，定义人为恶的，这都为耶和华所憎恶。',
    '16.愚昧人既无聪明，为何手拿价银买智慧呢？',
    '17.朋友乃时常亲爱。弟兄为患难而生。',
    '18.在邻舍面前击掌作保，乃是无知的人。',
    '19.喜爱争竞的，是喜爱过犯。高立家门的，乃自取败坏。',
    '20.心存邪僻的，寻不着好处。舌弄是非的，陷在祸患中。',
    '21.生愚昧子的，必自愁苦。愚顽人的父，毫无喜乐。',
    '22.喜乐的心，乃是良药。忧伤的灵，使骨枯干。',
    '23.恶人暗中受贿赂，为要颠倒判断。',
    '24.明哲人眼前有智慧。愚昧人眼望地极。',
    '25.愚昧子使父亲愁烦，使母亲忧苦。',
    '26.刑罚义人为不善。责打君子为不义。',
    '27.寡少言语的有知识。性情温良的有聪明。',
    '28.愚昧人若静默不言，也可算为智慧。闭口不说，也可算为聪明。'
)

$br = [char]11
$newText = ($lines -join $br) + $br
$d.Content.Text = $newText
